$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: merge the two runs "MON Mar 25" / " 14:19:48 IST 2019" into a
# single run of text "MON Mar 25 14:19:48 IST 2019".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("MON Mar 25 14:19:48 IST 2019", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "MON Mar 25 14:19:48 IST 2019", 2) | Out-Null

# ---------------------------------------------------------------------
# Edit 2: append a new "chick-in" purchase-details block (SAT MAR 30)
# right after the last "Amount Received mode ... - CASH" paragraph,
# i.e. right before the trailing blank paragraphs at the end of the
# document.
# ---------------------------------------------------------------------

# Locate the paragraph that ends the MON Mar 25 block: the last
# paragraph whose text is exactly "Amount Received mode<tabs>- CASH".
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd()
    if ($t -match "^Amount Received mode\t+- CASH$") {
        $targetIndex = $i
    }
}

$anchor = $d.Paragraphs.Item($targetIndex)
$anchor.Range.InsertParagraphAfter() | Out-Null

# --- paragraph: blank line --------------------------------------------------
$idx = $targetIndex + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter() | Out-Null

# --- paragraph: SAT MAR 30   15:12:29 IST 2019 -------------------------------
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("SAT MAR 30 15:12:29 IST 2019")
$p.Range.InsertParagraphAfter() | Out-Null

# --- paragraph: Person Name ... - TNR H --------------------------------------
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("Person Name`t`t`t`t- TNR H")
$p.Range.InsertParagraphAfter() | Out-Null

# --- paragraph: dashed separator ---------------------------------------------
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("---------------------------------------------------------------")
$p.Range.InsertParagraphAfter() | Out-Null

# --- paragraph: Item Name ... - CARROT ---------------------------------------
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("Item Name`t`t`t`t- CARROT")
$p.Range.InsertParagraphAfter() | Out-Null

# --- paragraph: Amount Received ... - 1564 (red text) ------------------------
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("Amount Received`t`t`t- 1564")
$p.Range.Font.Color = 255
$p.Range.InsertParagraphAfter() | Out-Null

# --- paragraph: Amount Received mode ... - CASH AND CLEARD -------------------
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("Amount Received mode`t`t- CASH AND CLEARD")
$p.Range.InsertParagraphAfter() | Out-Null

# --- paragraph: blank line -----------------------------------------------
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter() | Out-Null
